# Update BP terminal gate pricing (TGP) table: shift effective dates forward
# one day and refresh Diesel/ULP/PULP/e10 price columns (D:G) for rows 8-65.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8
$ws.Range("A8").Value = 46010
$ws.Range("D8").Value = 155.34
$ws.Range("E8").Value = 155.65
$ws.Range("F8").Value = 165.65
$ws.Range("G8").Value = 155.77

# Row 9
$ws.Range("A9").Value = 46010
$ws.Range("D9").Value = 155.34
$ws.Range("E9").Value = 155.65
$ws.Range("F9").Value = 165.65
$ws.Range("G9").Value = 155.77

# Row 10
$ws.Range("A10").Value = 46010
$ws.Range("D10").Value = 157.33
$ws.Range("E10").Value = 158.25
$ws.Range("F10").Value = 168.25
$ws.Range("G10").Value = 158.77

# Row 11
$ws.Range("A11").Value = 46009
$ws.Range("D11").Value = 156.56
$ws.Range("E11").Value = 156.89
$ws.Range("F11").Value = 166.89
$ws.Range("G11").Value = 157

# Row 12
$ws.Range("A12").Value = 46009
$ws.Range("D12").Value = 156.56
$ws.Range("E12").Value = 156.89
$ws.Range("F12").Value = 166.89
$ws.Range("G12").Value = 157

# Row 13
$ws.Range("A13").Value = 46009
$ws.Range("D13").Value = 158.63
$ws.Range("E13").Value = 159.46
$ws.Range("F13").Value = 169.46
$ws.Range("G13").Value = 159.98

# Row 17
$ws.Range("A17").Value = 46010
$ws.Range("D17").Value = 160.84
$ws.Range("E17").Value = 162.5
$ws.Range("F17").Value = 172.5

# Row 18
$ws.Range("A18").Value = 46009
$ws.Range("D18").Value = 161.69
$ws.Range("E18").Value = 163.15
$ws.Range("F18").Value = 173.15

# Row 22
$ws.Range("A22").Value = 46010
$ws.Range("D22").Value = 156.15
$ws.Range("E22").Value = 157.25
$ws.Range("F22").Value = 166.85
$ws.Range("G22").Value = 158.41

# Row 23
$ws.Range("A23").Value = 46010
$ws.Range("D23").Value = 162.68
$ws.Range("E23").Value = 162.57
$ws.Range("F23").Value = 172.57

# Row 24
$ws.Range("A24").Value = 46010
$ws.Range("D24").Value = 162.45
$ws.Range("E24").Value = 162.89
$ws.Range("F24").Value = 172.89

# Row 25
$ws.Range("A25").Value = 46010
$ws.Range("D25").Value = 162.95
$ws.Range("E25").Value = 162.53
$ws.Range("F25").Value = 172.53
$ws.Range("G25").Value = 162.3

# Row 26
$ws.Range("A26").Value = 46010
$ws.Range("D26").Value = 161.87
$ws.Range("E26").Value = 163.95
$ws.Range("F26").Value = 173.95

# Row 27
$ws.Range("A27").Value = 46009
$ws.Range("D27").Value = 157.37
$ws.Range("E27").Value = 158.46
$ws.Range("F27").Value = 168.06
$ws.Range("G27").Value = 159.62

# Row 28
$ws.Range("A28").Value = 46009
$ws.Range("D28").Value = 163.55
$ws.Range("E28").Value = 163.57
$ws.Range("F28").Value = 173.57

# Row 29
$ws.Range("A29").Value = 46009
$ws.Range("D29").Value = 163.32
$ws.Range("E29").Value = 163.89
$ws.Range("F29").Value = 173.89

# Row 30
$ws.Range("A30").Value = 46009
$ws.Range("D30").Value = 163.82
$ws.Range("E30").Value = 163.53
$ws.Range("F30").Value = 173.53
$ws.Range("G30").Value = 163.3

# Row 31
$ws.Range("A31").Value = 46009
$ws.Range("D31").Value = 162.73
$ws.Range("E31").Value = 164.96
$ws.Range("F31").Value = 174.96

# Row 35
$ws.Range("A35").Value = 46010
$ws.Range("D35").Value = 155.84
$ws.Range("E35").Value = 155.6
$ws.Range("F35").Value = 164.6

# Row 36
$ws.Range("A36").Value = 46009
$ws.Range("D36").Value = 156.92
$ws.Range("E36").Value = 156.6
$ws.Range("F36").Value = 165.6

# Row 40
$ws.Range("A40").Value = 46010
$ws.Range("D40").Value = 162.15
$ws.Range("E40").Value = 163.13
$ws.Range("F40").Value = 173.13

# Row 41
$ws.Range("A41").Value = 46010
$ws.Range("D41").Value = 161.86
$ws.Range("E41").Value = 163.55
$ws.Range("F41").Value = 173.55

# Row 42
$ws.Range("A42").Value = 46009
$ws.Range("D42").Value = 163.02
$ws.Range("E42").Value = 163.81
$ws.Range("F42").Value = 173.81

# Row 43
$ws.Range("A43").Value = 46009
$ws.Range("D43").Value = 162.74
$ws.Range("E43").Value = 164.23
$ws.Range("F43").Value = 174.23

# Row 47
$ws.Range("A47").Value = 46010
$ws.Range("D47").Value = 157.06
$ws.Range("E47").Value = 157.77
$ws.Range("F47").Value = 167.77

# Row 48
$ws.Range("A48").Value = 46010
$ws.Range("D48").Value = 156.87
$ws.Range("E48").Value = 157.86
$ws.Range("F48").Value = 167.86

# Row 49
$ws.Range("A49").Value = 46009
$ws.Range("D49").Value = 158.11
$ws.Range("E49").Value = 158.53
$ws.Range("F49").Value = 168.53

# Row 50
$ws.Range("A50").Value = 46009
$ws.Range("D50").Value = 157.92
$ws.Range("E50").Value = 158.62
$ws.Range("F50").Value = 168.62

# Row 54
$ws.Range("A54").Value = 46010
$ws.Range("D54").Value = 171.66
$ws.Range("E54").Value = 173.28
$ws.Range("F54").Value = 183.28

# Row 55
$ws.Range("A55").Value = 46010
$ws.Range("D55").Value = 159.84
$ws.Range("E55").Value = 168.05
$ws.Range("F55").Value = 178.05

# Row 56
$ws.Range("A56").Value = 46010
$ws.Range("D56").Value = 161.89

# Row 57
$ws.Range("A57").Value = 46010
$ws.Range("D57").Value = 160.93
$ws.Range("E57").Value = 162.33

# Row 58
$ws.Range("A58").Value = 46010
$ws.Range("D58").Value = 156.83
$ws.Range("E58").Value = 158.37
$ws.Range("F58").Value = 168.37

# Row 59
$ws.Range("A59").Value = 46010
$ws.Range("D59").Value = 163.76
$ws.Range("E59").Value = 170.56

# Row 60
$ws.Range("A60").Value = 46009
$ws.Range("D60").Value = 172.51
$ws.Range("E60").Value = 174.33
$ws.Range("F60").Value = 184.33

# Row 61
$ws.Range("A61").Value = 46009
$ws.Range("D61").Value = 160.71
$ws.Range("E61").Value = 169.25
$ws.Range("F61").Value = 179.25

# Row 62
$ws.Range("A62").Value = 46009
$ws.Range("D62").Value = 162.97

# Row 63
$ws.Range("A63").Value = 46009
$ws.Range("D63").Value = 161.99
$ws.Range("E63").Value = 163.52

# Row 64
$ws.Range("A64").Value = 46009
$ws.Range("D64").Value = 157.89
$ws.Range("E64").Value = 159.57
$ws.Range("F64").Value = 169.57

# Row 65
$ws.Range("A65").Value = 46009
$ws.Range("D65").Value = 164.61
$ws.Range("E65").Value = 171.59
